$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 2-10 (columns B:G) down by one row to rows 3-11,
# starting from the bottom so we don't overwrite values before they're copied.
for ($r = 10; $r -ge 2; $r--) {
    $srcRow = $r
    $dstRow = $r + 1
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($dstRow, $c).Value = $ws.Cells.Item($srcRow, $c).Value2
    }
}

# Set the new values for row 2 (columns B:G)
$ws.Cells.Item(2, 2).Value = 0.001217196552899444
$ws.Cells.Item(2, 3).Value = 2.134347703298304
$ws.Cells.Item(2, 4).Value = 20.89521853204074
$ws.Cells.Item(2, 5).Value = 4.571128802827672
$ws.Cells.Item(2, 6).Value = 4.673863460101999
$ws.Cells.Item(2, 7).Value = 23
